$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.135806
$ws.Range("H2").Value = 0.407418
$ws.Range("M2").Value = 1.701929666666667
$ws.Range("N2").Value = 5.105789
$ws.Range("O2").Value = 0.02105622887134972
$ws.Range("P2").Value = 0.02105622887134972
$ws.Range("Q2").Value = 0.2311322603113333
$ws.Range("R2").Value = 2.080190342802
$ws.Range("S2").Value = 0.02105622887134972
$ws.Range("T2").Value = 0.02105622887134972

$ws.Range("G3").Value = 0.135806
$ws.Range("H3").Value = 0.407418
$ws.Range("O3").Value = 0.7732971809418951
$ws.Range("P3").Value = 0.7732971809418953
$ws.Range("Q3").Value = 8.488411026282002
$ws.Range("R3").Value = 76.39569923653801
$ws.Range("S3").Value = 0.7732971809418951
$ws.Range("T3").Value = 0.7732971809418953

$ws.Range("G4").Value = 0.135806
$ws.Range("H4").Value = 0.407418
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4338690000000001
$ws.Range("N4").Value = 1.301607
$ws.Range("O4").Value = 0.005367815805265532
$ws.Range("P4").Value = 0.005367815805265533
$ws.Range("Q4").Value = 0.05892201341400001
$ws.Range("R4").Value = 0.530298120726
$ws.Range("S4").Value = 0.005367815805265532
$ws.Range("T4").Value = 0.005367815805265533

$ws.Range("G5").Value = 0.135806
$ws.Range("H5").Value = 0.407418
$ws.Range("M5").Value = 15.972384
$ws.Range("N5").Value = 47.917152
$ws.Range("O5").Value = 0.1976099128607259
$ws.Range("P5").Value = 0.1976099128607259
$ws.Range("Q5").Value = 2.169145581504
$ws.Range("R5").Value = 19.522310233536
$ws.Range("S5").Value = 0.1976099128607259
$ws.Range("T5").Value = 0.1976099128607259

$ws.Range("G6").Value = 0.135806
$ws.Range("H6").Value = 0.407418
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2157183333333333
$ws.Range("N6").Value = 0.647155
$ws.Range("O6").Value = 0.002668861520763652
$ws.Range("P6").Value = 0.002668861520763652
$ws.Range("Q6").Value = 0.02929584397666667
$ws.Range("R6").Value = 0.26366259579
$ws.Range("S6").Value = 0.002668861520763652
$ws.Range("T6").Value = 0.002668861520763652
